# Update cryptocurrency price and volume(1h) data in the active worksheet.
# For the Price column (D), we force a Text number format before assigning
# the value so that values which look numeric (e.g. "1.002") are preserved
# exactly as text, matching how the source data was originally stored.
# We then reset the cell style back to Normal so no extra formatting is
# left behind, while the cell keeps its text (string) value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.723.39"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +7.03%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.947.17"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +5.39%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.68%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "341.63"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.78%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.61%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4781"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.80%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4124"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +6.97%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "47.71"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.90%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08237"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.37%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.034"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +6.97%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.65"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +6.56%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.939.54"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.85%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.137"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.20%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.353"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.99%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "91.64"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.72%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.003"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.69%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001057"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.80%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06667"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.00%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.01"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.30%  "
$ws.Range("E21").Value = "  -0.70%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "29.698.85"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +6.94%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.582"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.29%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.23"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.82%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.289"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.51%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.177.65"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.30%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "161.06"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.43%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.17"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.87%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.168"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.98%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.639"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.88%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "122.96"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.62%  "
$ws.Range("E32").Value = "  +7.12%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09656"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.33%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.470"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +11.00%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.680"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.20%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.488"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.45%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06255"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.86%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02316"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.60%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.482"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.17%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.187"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.13%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6070"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.51%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "10.70"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +6.41%  "
$ws.Range("E43").Value = "  -0.64%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1896"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.71%  "
$ws.Range("E45").Value = "  -1.97%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.374"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +32.40%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5710"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.82%  "
$ws.Range("E48").Value = "  +4.08%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.07416"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +8.34%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.987"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.63%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "112.76"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.85%  "

